$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 11.5
$ws.Range("L2").Value = 1.43
$ws.Range("V2").Value = 1.09
$ws.Range("W2").Value = 3.6
$ws.Range("X2").Value = 14.5
$ws.Range("Y2").Value = 28
$ws.Range("Z2").Value = 110
$ws.Range("AA2").Value = 750
$ws.Range("AB2").Value = 6.2
$ws.Range("AC2").Value = 12
$ws.Range("AD2").Value = 48
$ws.Range("AE2").Value = 320
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 44
$ws.Range("AI2").Value = 270
$ws.Range("AJ2").Value = 10
$ws.Range("AK2").Value = 17.5
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 360
$ws.Range("AN2").Value = 8
$ws.Range("AO2").Value = 620
$ws.Range("L3").Value = 1.4
$ws.Range("M3").Value = 1.08
$ws.Range("T3").Value = 1.79
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 1.52
$ws.Range("X3").Value = 13
$ws.Range("AB3").Value = 11.5
$ws.Range("AD3").Value = 12
$ws.Range("AF3").Value = 18.5
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 42
$ws.Range("AJ3").Value = 44
$ws.Range("AK3").Value = 32
$ws.Range("AL3").Value = 44
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 28
$ws.Range("AO3").Value = 25
$ws.Range("S4").Value = 2.16
$ws.Range("AC4").Value = 13
$ws.Range("AM4").Value = 85
$ws.Range("T5").Value = 1.67
$ws.Range("U5").Value = 2.38
$ws.Range("X5").Value = 19
$ws.Range("Y5").Value = 16.5
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 75
$ws.Range("AB5").Value = 12
$ws.Range("AD5").Value = 15
$ws.Range("AE5").Value = 38
$ws.Range("AF5").Value = 15.5
$ws.Range("AI5").Value = 42
$ws.Range("AJ5").Value = 28
$ws.Range("AK5").Value = 22
$ws.Range("AL5").Value = 32
$ws.Range("AN5").Value = 13.5
$ws.Range("AO5").Value = 30
$ws.Range("F6").Value = 1.8
$ws.Range("I6").Value = 5.6
$ws.Range("J6").Value = 3.7
$ws.Range("P6").Value = 1.8
$ws.Range("X6").Value = 12.5
$ws.Range("Y6").Value = 16.5
$ws.Range("Z6").Value = 42
$ws.Range("AA6").Value = 190
$ws.Range("AD6").Value = 22
$ws.Range("AE6").Value = 85
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 110
$ws.Range("AJ6").Value = 19.5
$ws.Range("AK6").Value = 19.5
$ws.Range("AL6").Value = 44
$ws.Range("AM6").Value = 180
$ws.Range("AN6").Value = 14.5
$ws.Range("AO6").Value = 150
$ws.Range("H7").Value = 2.42
$ws.Range("I7").Value = 2.44
$ws.Range("S7").Value = 3.6
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 16
$ws.Range("AA7").Value = 34
$ws.Range("AB7").Value = 13.5
$ws.Range("AE7").Value = 25
$ws.Range("AF7").Value = 23
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 40
$ws.Range("AJ7").Value = 55
$ws.Range("AL7").Value = 48
$ws.Range("AN7").Value = 34
$ws.Range("AO7").Value = 21
$ws.Range("N8").Value = 3.55
$ws.Range("O8").Value = 1.37
$ws.Range("Q8").Value = 2.12
$ws.Range("U8").Value = 2
$ws.Range("X8").Value = 12.5
$ws.Range("Y8").Value = 15
$ws.Range("Z8").Value = 34
$ws.Range("AA8").Value = 130
$ws.Range("AD8").Value = 18.5
$ws.Range("AE8").Value = 65
$ws.Range("AF8").Value = 11.5
$ws.Range("AH8").Value = 22
$ws.Range("AI8").Value = 70
$ws.Range("AL8").Value = 40
$ws.Range("AM8").Value = 100
$ws.Range("AN8").Value = 15.5
$ws.Range("AO8").Value = 100
$ws.Range("F9").Value = 3.6
$ws.Range("G9").Value = 3.65
$ws.Range("X9").Value = 17
$ws.Range("AB9").Value = 15.5
$ws.Range("AE9").Value = 20
$ws.Range("AF9").Value = 27
$ws.Range("AG9").Value = 15
$ws.Range("AJ9").Value = 65
$ws.Range("AK9").Value = 40
$ws.Range("AL9").Value = 46
$ws.Range("AN9").Value = 36
$ws.Range("AO9").Value = 15
$ws.Range("F10").Value = 2.52
$ws.Range("G10").Value = 2.54
$ws.Range("H10").Value = 3.35
$ws.Range("O10").Value = 1.45
$ws.Range("Q10").Value = 2.38
$ws.Range("X10").Value = 10.5
$ws.Range("Y10").Value = 11
$ws.Range("Z10").Value = 22
$ws.Range("AA10").Value = 70
$ws.Range("AB10").Value = 8.800000000000001
$ws.Range("AD10").Value = 15
$ws.Range("AE10").Value = 48
$ws.Range("AF10").Value = 15
$ws.Range("AG10").Value = 12.5
$ws.Range("AH10").Value = 21
$ws.Range("AI10").Value = 70
$ws.Range("AJ10").Value = 36
$ws.Range("AK10").Value = 32
$ws.Range("AL10").Value = 55
$ws.Range("AM10").Value = 170
$ws.Range("AN10").Value = 32
$ws.Range("AO10").Value = 60
$ws.Range("N11").Value = 3
$ws.Range("P11").Value = 1.67
$ws.Range("R11").Value = 1.24
$ws.Range("S11").Value = 4.7
$ws.Range("U11").Value = 1.92
$ws.Range("Y11").Value = 11.5
$ws.Range("Z11").Value = 24
$ws.Range("AA11").Value = 75
$ws.Range("AB11").Value = 8.199999999999999
$ws.Range("AD11").Value = 16
$ws.Range("AE11").Value = 55
$ws.Range("AF11").Value = 14.5
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 70
$ws.Range("AJ11").Value = 36
$ws.Range("AK11").Value = 32
$ws.Range("AL11").Value = 55
$ws.Range("AM11").Value = 170
$ws.Range("AN11").Value = 30
$ws.Range("AO11").Value = 65
$ws.Range("P12").Value = 2.38
$ws.Range("Q12").Value = 1.69
$ws.Range("R12").Value = 1.55
$ws.Range("S12").Value = 2.72
$ws.Range("T12").Value = 1.97
$ws.Range("X12").Value = 23
$ws.Range("Y12").Value = 9.4
$ws.Range("Z12").Value = 8.800000000000001
$ws.Range("AA12").Value = 12
$ws.Range("AB12").Value = 32
$ws.Range("AC12").Value = 12
$ws.Range("AE12").Value = 15
$ws.Range("AG12").Value = 34
$ws.Range("AH12").Value = 26
$ws.Range("AI12").Value = 36
$ws.Range("AJ12").Value = 380
$ws.Range("AK12").Value = 170
$ws.Range("AM12").Value = 160
$ws.Range("AN12").Value = 210
$ws.Range("AO12").Value = 5.7
$ws.Range("F13").Value = 1.72
$ws.Range("G13").Value = 1.73
$ws.Range("I13").Value = 5.7
$ws.Range("P13").Value = 2.14
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 3.2
$ws.Range("U13").Value = 2.12
$ws.Range("Y13").Value = 21
$ws.Range("AE13").Value = 75
$ws.Range("AK13").Value = 17.5
$ws.Range("AO13").Value = 100
$ws.Range("F14").Value = 5.4
$ws.Range("G14").Value = 5.5
$ws.Range("I14").Value = 1.76
$ws.Range("K14").Value = 4.2
$ws.Range("P14").Value = 2.18
$ws.Range("X14").Value = 19
$ws.Range("AA14").Value = 18
$ws.Range("AC14").Value = 9
$ws.Range("AI14").Value = 34
$ws.Range("AK14").Value = 70
$ws.Range("AN14").Value = 75
$ws.Range("AO14").Value = 9.199999999999999
